# Apply updated power-flow line results (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.3145716436525561
    $ws.Range("C2").Value = 0.1665262621700663
    $ws.Range("E2").Value = 0.1274432258134865
    $ws.Range("F2").Value = 0.4443680307746263
    $ws.Range("G2").Value = 0.002532364749321041
    $ws.Range("I2").Value = 2.090370996350273
    $ws.Range("K2").Value = 0.4804745443254319
    $ws.Range("M2").Value = 0.2750764028036983
    $ws.Range("B3").Value = 0.2960428588898765
    $ws.Range("C3").Value = 0.153672118406206
    $ws.Range("E3").Value = 0.1164969290119231
    $ws.Range("F3").Value = 0.387822817061874
    $ws.Range("G3").Value = 0.002536723808828535
    $ws.Range("I3").Value = 2.00773936403391
    $ws.Range("K3").Value = 0.4490860054193035
    $ws.Range("M3").Value = 0.2543205311393919
    $ws.Range("B4").Value = 0.284972551378587
    $ws.Range("C4").Value = 0.145851506622364
    $ws.Range("E4").Value = 0.1098447830212308
    $ws.Range("F4").Value = 0.3531389305168915
    $ws.Range("G4").Value = 0.002539536973484404
    $ws.Range("I4").Value = 1.957149556275454
    $ws.Range("K4").Value = 0.4301845372416153
    $ws.Range("M4").Value = 0.24175707665065
    $ws.Range("B5").Value = 0.2805380868041567
    $ws.Range("C5").Value = 0.1426823597828957
    $ws.Range("E5").Value = 0.1071510078412103
    $ws.Range("F5").Value = 0.3390132514313251
    $ws.Range("G5").Value = 0.002540717850618583
    $ws.Range("I5").Value = 1.936568675502571
    $ws.Range("K5").Value = 0.4225748587425926
    $ws.Range("M5").Value = 0.2366823460395153
    $ws.Range("B6").Value = 0.2798063763581951
    $ws.Range("C6").Value = 0.142157193866808
    $ws.Range("E6").Value = 0.1067047281865001
    $ws.Range("F6").Value = 0.336668177824194
    $ws.Range("G6").Value = 0.002540916020977355
    $ws.Range("I6").Value = 1.933153286165066
    $ws.Range("K6").Value = 0.4213168686719655
    $ws.Range("M6").Value = 0.2358423939677934
    $ws.Range("B7").Value = 0.2849124361451629
    $ws.Range("C7").Value = 0.1458086946041703
    $ws.Range("E7").Value = 0.109808385276466
    $ws.Range("F7").Value = 0.3529483938344953
    $ws.Range("G7").Value = 0.002539552759336142
    $ws.Range("I7").Value = 1.956871856662573
    $ws.Range("K7").Value = 0.4300815353039127
    $ws.Range("M7").Value = 0.2416884555524135
    $ws.Range("B8").Value = 0.3081191378893493
    $ws.Range("C8").Value = 0.1620790726179848
    $ws.Range("E8").Value = 0.1236544507638584
    $ws.Range("F8").Value = 0.4248636149813336
    $ws.Range("G8").Value = 0.002533839456257255
    $ws.Range("I8").Value = 2.06184805438977
    $ws.Range("K8").Value = 0.469574399872414
    $ws.Range("M8").Value = 0.2678819765889671
    $ws.Range("B9").Value = 0.3560749279322408
    $ws.Range("C9").Value = 0.1945690763870118
    $ws.Range("E9").Value = 0.1513685845059172
    $ws.Range("F9").Value = 0.5661985755041457
    $ws.Range("G9").Value = 0.002523714749544259
    $ws.Range("I9").Value = 2.268963317160797
    $ws.Range("K9").Value = 0.5499939626989487
    $ws.Range("M9").Value = 0.3207056048756911
    $ws.Range("B10").Value = 0.3928271699457468
    $ws.Range("C10").Value = 0.2188165761376126
    $ws.Range("E10").Value = 0.1720959088970417
    $ws.Range("F10").Value = 0.6702781546542269
    $ws.Range("G10").Value = 0.002516926268259536
    $ws.Range("I10").Value = 2.422043684258085
    $ws.Range("K10").Value = 0.6109384354423071
    $ws.Range("M10").Value = 0.3604419534209882
    $ws.Range("B11").Value = 0.4098828557474405
    $ws.Range("C11").Value = 0.229934030187195
    $ws.Range("E11").Value = 0.1816098812060716
    $ws.Range("F11").Value = 0.7176906081379002
    $ws.Range("G11").Value = 0.002513977540907986
    $ws.Range("I11").Value = 2.491914380318462
    $ws.Range("K11").Value = 0.6390783707849437
    $ws.Range("M11").Value = 0.3787288320221762
    $ws.Range("B12").Value = 0.4163902788602911
    $ws.Range("C12").Value = 0.2341567894197567
    $ws.Range("E12").Value = 0.1852251783181771
    $ws.Range("F12").Value = 0.7356546913071611
    $ws.Range("G12").Value = 0.002512880853209496
    $ws.Range("I12").Value = 2.51840849873787
    $ws.Range("K12").Value = 0.6497948003174372
    $ws.Range("M12").Value = 0.3856844972706952
    $ws.Range("B13").Value = 0.4149866124205062
    $ws.Range("C13").Value = 0.2332467678121475
    $ws.Range("E13").Value = 0.1844459952581801
    $ws.Range("F13").Value = 0.7317853510981394
    $ws.Range("G13").Value = 0.002513116159819103
    $ws.Range("I13").Value = 2.512700900956844
    $ws.Range("K13").Value = 0.6474841269334206
    $ws.Range("M13").Value = 0.3841850909106128
    $ws.Range("B14").Value = 0.4104172451465331
    $ws.Range("C14").Value = 0.2302811805135434
    $ws.Range("E14").Value = 0.1819070601558224
    $ws.Range("F14").Value = 0.7191683204515869
    $ws.Range("G14").Value = 0.002513886916939573
    $ws.Range("I14").Value = 2.494093342385355
    $ws.Range("K14").Value = 0.6399588027993559
    $ws.Range("M14").Value = 0.3793004579884851
    $ws.Range("B15").Value = 0.4076247424315227
    $ws.Range("C15").Value = 0.2284663509483664
    $ws.Range("E15").Value = 0.1803535351861854
    $ws.Range("F15").Value = 0.7114413442032514
    $ws.Range("G15").Value = 0.002514361620451328
    $ws.Range("I15").Value = 2.48270037140918
    $ws.Range("K15").Value = 0.6353572153183791
    $ws.Range("M15").Value = 0.3763125091470911
    $ws.Range("B16").Value = 0.3917193505712078
    $ws.Range("C16").Value = 0.2180918025904077
    $ws.Range("E16").Value = 0.1714758878288265
    $ws.Range("F16").Value = 0.6671810134426437
    $ws.Range("G16").Value = 0.002517121769641372
    $ws.Range("I16").Value = 2.417482342895369
    $ws.Range("K16").Value = 0.6091078486015533
    $ws.Range("M16").Value = 0.359251146538341
    $ws.Range("B17").Value = 0.3820484492137268
    $ws.Range("C17").Value = 0.2117498777247988
    $ws.Range("E17").Value = 0.1660517506855328
    $ws.Range("F17").Value = 0.6400460337215605
    $ws.Range("G17").Value = 0.002518850652008747
    $ws.Range("I17").Value = 2.377534465191843
    $ws.Range("K17").Value = 0.5931116993140506
    $ws.Range("M17").Value = 0.3488388745385649
    $ws.Range("B18").Value = 0.376517691694346
    $ws.Range("C18").Value = 0.208110355222118
    $ws.Range("E18").Value = 0.1629399113835674
    $ws.Range("F18").Value = 0.6244449056556647
    $ws.Range("G18").Value = 0.002519858186046349
    $ws.Range("I18").Value = 2.354579376673712
    $ws.Range("K18").Value = 0.5839502537322119
    $ws.Range("M18").Value = 0.3428698025638468
    $ws.Range("B19").Value = 0.3746505041296189
    $ws.Range("C19").Value = 0.2068794724586667
    $ws.Range("E19").Value = 0.1618876565633229
    $ws.Range("F19").Value = 0.6191636801734006
    $ws.Range("G19").Value = 0.002520201577447795
    $ws.Range("I19").Value = 2.346810874147138
    $ws.Range("K19").Value = 0.5808550450563246
    $ws.Range("M19").Value = 0.340852160402946
    $ws.Range("B20").Value = 0.3830746509093785
    $ws.Range("C20").Value = 0.2124241372551694
    $ws.Range("E20").Value = 0.1666283309074359
    $ws.Range("F20").Value = 0.642933953830422
    $ws.Range("G20").Value = 0.002518665251713153
    $ws.Range("I20").Value = 2.381784706607078
    $ws.Range("K20").Value = 0.5948104638121663
    $ws.Range("M20").Value = 0.3499452264707017
    $ws.Range("B21").Value = 0.411758051613333
    $ws.Range("C21").Value = 0.2311518951157723
    $ws.Range("E21").Value = 0.1826524632869919
    $ws.Range("F21").Value = 0.7228739723492197
    $ws.Range("G21").Value = 0.002513659987029533
    $ws.Range("I21").Value = 2.499557851185074
    $ws.Range("K21").Value = 0.6421675273289793
    $ws.Range("M21").Value = 0.3807343522958817
    $ws.Range("B22").Value = 0.4307889593627294
    $ws.Range("C22").Value = 0.2434664332441798
    $ws.Range("E22").Value = 0.1931985427579619
    $ws.Range("F22").Value = 0.7751780083420101
    $ws.Range("G22").Value = 0.002510504878991324
    $ws.Range("I22").Value = 2.576737858826675
    $ws.Range("K22").Value = 0.6734707377963218
    $ws.Range("M22").Value = 0.4010367301351678
    $ws.Range("B23").Value = 0.4206056392731057
    $ws.Range("C23").Value = 0.2368869856308322
    $ws.Range("E23").Value = 0.1875630742808312
    $ws.Range("F23").Value = 0.7472568307915566
    $ws.Range("G23").Value = 0.002512178231351757
    $ws.Range("I23").Value = 2.53552570057002
    $ws.Range("K23").Value = 0.6567311480501417
    $ws.Range("M23").Value = 0.3901843286091804
    $ws.Range("B24").Value = 0.3826106142814467
    $ws.Range("C24").Value = 0.2121192841914592
    $ws.Range("E24").Value = 0.1663676385555561
    $ws.Range("F24").Value = 0.6416283278902313
    $ws.Range("G24").Value = 0.002518749028863915
    $ws.Range("I24").Value = 2.379863136857438
    $ws.Range("K24").Value = 0.5940423435560263
    $ws.Range("M24").Value = 0.3494449915993414
    $ws.Range("B25").Value = 0.3428368345834656
    $ws.Range("C25").Value = 0.1857149904068365
    $ws.Range("E25").Value = 0.1438086317641307
    $ws.Range("F25").Value = 0.5279251897347308
    $ws.Range("G25").Value = 0.002526339026249512
    $ws.Range("I25").Value = 2.212783439988314
    $ws.Range("K25").Value = 0.5279151188793207
    $ws.Range("M25").Value = 0.3062554185808253

